$d = $word.ActiveDocument

$old0 = "Aircraft certification requires extensive testing including flyover noise measurements. The measurements are acquired in remote locations to minimize environmental noise contamination from altering the recorded noise levels that are submitting to regulatory agencies such as the FAA or causing costly repeat flyovers to acquire clean recordings for the certification process. Typical sources of contamination include bird chips and other wildlife or livestock vocalizations, insect noise, and traffic noise. The currently process involves using three engineers to listen to the live microphone feeds and alert the test manager of any noise contamination. The goal of this project is to study the feasibility of machine learning algorithms to perform the identification task. The will include a wide survey of feature sets and classification algorithms as to rank the combinations. "
$new0 = "Aircraft certification requires extensive testing including flyover noise measurements. The measurements are acquired in remote locations to minimize contamination from environmental noise that may alter the recorded noise levels that are then submitted to regulatory agencies such as the FAA. Contamination may also cause costly repeat flyovers to acquire clean recordings for the certification process. Typical sources of contamination include bird chirps, other wildlife or livestock vocalizations, insect noise, and traffic noise. The current process involves using three engineers to listen to the live microphone feeds and alert the test manager of any noise contamination. The goal of this project is to study the feasibility of machine learning algorithms performing the task of identifying noise contamination. To understand this, our project will include a wide survey of feature sets and classification algorithms to rank the combinations and effectiveness of each. "
$found0 = $d.Content.Find.Execute($old0, $true, $false, $false, $false, $false, $true, 1, $false, $new0, 2)
Write-Output "Replace 0: $found0"

$old1 = "all data is expected to be in file format instead of streams or API calls. Extra data may be needed to expand the"
$new1 = "all data is expected to be in file format instead of streams or API calls. Additional sample data may be required to expand the"
$found1 = $d.Content.Find.Execute($old1, $true, $false, $false, $false, $false, $true, 1, $false, $new1, 2)
Write-Output "Replace 1: $found1"

$old2 = "The data were recorded on three separate days with one day each in the months of November 2017"
$new2 = "The data were recorded on three separate days in the months of November 2017"
$found2 = $d.Content.Find.Execute($old2, $true, $false, $false, $false, $false, $true, 1, $false, $new2, 2)
Write-Output "Replace 2: $found2"

$old3 = "This should not impact the project as the levels are to be normalized to remove the overall signal energy to remove the effects of differences in propagation distances from the sources to the microphone sensors."
$new3 = "This should not impact the project as the levels are to be normalized to remove the effects of differences in propagation distances from sound sources to the microphone sensors."
$found3 = $d.Content.Find.Execute($old3, $true, $false, $false, $false, $false, $true, 1, $false, $new3, 2)
Write-Output "Replace 3: $found3"

$old4 = ", respectively. "
$new4 = " respectively. "
$found4 = $d.Content.Find.Execute($old4, $true, $false, $false, $false, $false, $true, 1, $false, $new4, 2)
Write-Output "Replace 4: $found4"

$old5 = "As"
$new5 = "Despite"
$found5 = $d.Content.Find.Execute($old5, $true, $true, $false, $false, $false, $true, 1, $false, $new5, 2)
Write-Output "Replace 5: $found5"

$old6 = "the data is of high quality, there are still concerns with its use in this project. The recorded data is heavily skewed towards aircraft signatures thus creating an imbalance in the number of samples for each class. Also, each signal can produce 50-100 blocks or samples of features for classification and each block may contain 10’s to 1,000’s of feature depending on the feature generation option. This will present a data management concern as"
$new6 = "the data being of high quality, there are concerns with its use in this project. The recorded data is heavily skewed towards aircraft signatures thus creating an imbalance in the number of samples for each class. Additionally, each signal can produce 50-100 blocks or samples of features for classification and each block may contain 10’s to 1,000’s of features depending on the how features are generated. This will present a data management concern as"
$found6 = $d.Content.Find.Execute($old6, $true, $false, $false, $false, $false, $true, 1, $false, $new6, 2)
Write-Output "Replace 6: $found6"

$old7 = "the data, features, and classification labels must be accuracy tracked through the data processing."
$new7 = "the data, features, and classification labels must be accurately tracked through the data processing."
$found7 = $d.Content.Find.Execute($old7, $true, $false, $false, $false, $false, $true, 1, $false, $new7, 2)
Write-Output "Replace 7: $found7"

$old8 = "the impact the signal-to-noise ratio has on the system performance. Also, public domain recordings are available for use to extend the set of wildlife/livestock vocalizations"
$new8 = "the impact that the signal-to-noise ratio has on the system performance. If these combinations aren’t sufficient, public domain recordings are available for use to extend the set of wildlife/livestock vocalizations"
$found8 = $d.Content.Find.Execute($old8, $true, $false, $false, $false, $false, $true, 1, $false, $new8, 2)
Write-Output "Replace 8: $found8"

$old9 = "overlapped by 25% for octave analysis resulting in 79 spectral estimates."
$new9 = "overlapped by 25% of the record width for octave analysis resulting in 79 spectral estimates."
$found9 = $d.Content.Find.Execute($old9, $true, $false, $false, $false, $false, $true, 1, $false, $new9, 2)
Write-Output "Replace 9: $found9"

$old10 = "The octave spectrum presented in the examples above demonstrate one of many different feature set that could be used for classification between the data sets. Other feature generation options include 1/n octave spectrum, which would break the octave bands seen above into smaller bands, cepstrum processing, short-time Fourier transforms, and wavelets. All the different feature set have the potential to provide varying degrees of separation between the signals. This project should survey as many feature sets as practical to gain the understanding of which provide the best separation. Although, the classification algorithm will also impact the system performance and should be studied as well. The various combination of the feature sets and classification algorithm will form a rich test matrix for this project to study. "
$new10 = "The octave spectrum presented in the examples above demonstrate one of many different feature sets that could be used for classification of environmental noise. Other feature generation options include 1/n octave spectrum, which would break the octave bands seen above into smaller bands, cepstrum processing, short-time Fourier transforms, and wavelets. All the different feature sets have the potential to provide varying degrees of separation between the signals. This project should survey as many feature sets as practical to gain the understanding of which provide the best separation. Selection of the classification algorithm will also impact the system performance and should be studied as well. The various combinations of a feature set and classification algorithm will form a rich test matrix for this project to study. "
$found10 = $d.Content.Find.Execute($old10, $true, $false, $false, $false, $false, $true, 1, $false, $new10, 2)
Write-Output "Replace 10: $found10"

$old11 = "The raw data need for this project is represented as audio time series recordings and is readily available. The data set provided by the sponsor is representative but may need to be augmented to address the class imbalance. The data issues are present are mostly confined to data management practices as numerous different combinations of signals, processing techniques, and classifiers will need be tracked and the results reported on. However, the overall goal of the project supports the learning opportunity present with such a diverse set of features and algorithms. "
$new11 = "The raw data needed for this project is represented as audio time series recordings and is readily available. The data set provided by the sponsor is representative but may need to be augmented to address the class imbalance. The data issues that are present are mostly confined to data management practices as numerous combinations of signals, processing techniques, and classifiers will need be tracked and the results reported on. With such a diverse set of features and algorithms to explore, this project presents a great learning opportunity for those interested in signal processing."
$found11 = $d.Content.Find.Execute($old11, $true, $false, $false, $false, $false, $true, 1, $false, $new11, 2)
Write-Output "Replace 11: $found11"

if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

$gbAnchor = $d.Content
$gbAnchor.Find.Execute("classification algorithm will also impact the system performance a") | Out-Null
$gbRange = $d.Range($gbAnchor.End, $gbAnchor.End)
$d.Bookmarks.Add("_GoBack", $gbRange) | Out-Null
Write-Output "GoBack re-added at $($gbRange.Start)"

